$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 185
$ws.Range("I4").Value = 185
$ws.Range("K4").Value = 185
$ws.Range("M4").Value = -71
$ws.Range("H33").Value = 128.58333
$ws.Range("I33").Value = 85.14286
$ws.Range("J33").Value = 189.4
$ws.Range("K33").Value = 85.14286
$ws.Range("L33").Value = 189.4
$ws.Range("M33").Value = 143.85714
$ws.Range("N33").Value = -647.4
$ws.Range("H48").Value = 3980.5
$ws.Range("I48").Value = 3922
$ws.Range("J48").Value = 4000
$ws.Range("K48").Value = 11766
$ws.Range("L48").Value = 12000
$ws.Range("M48").Value = -11474
$ws.Range("N48").Value = -12584
$ws.Range("H56").Value = 3980.5
$ws.Range("I56").Value = 3922
$ws.Range("J56").Value = 4000
$ws.Range("K56").Value = 11766
$ws.Range("L56").Value = 12000
$ws.Range("M56").Value = -11232
$ws.Range("N56").Value = -13068
$ws.Range("H116").Value = 2680
$ws.Range("I116").Value = 2516
$ws.Range("K116").Value = 2516
$ws.Range("M116").Value = 926
$ws.Range("N116").ClearContents()
$ws.Range("H129").Value = 1036.625
$ws.Range("I129").Value = 1041.8572
$ws.Range("J129").Value = 1000
$ws.Range("K129").Value = 3125.5716
$ws.Range("L129").Value = 3000
$ws.Range("M129").Value = 1874.4284
$ws.Range("N129").Value = -13000
$ws.Range("H137").Value = 4952.1904
$ws.Range("I137").Value = 1999.625
$ws.Range("J137").Value = 6769.154
$ws.Range("K137").Value = 5998.875
$ws.Range("L137").Value = 20307.462
$ws.Range("M137").Value = -3448.875
$ws.Range("N137").Value = -25407.462

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 740.0769
$ws.Range("I2").Value = 685.0833
$ws.Range("K2").Value = 685.0833
$ws.Range("M2").Value = -572.0833
$ws.Range("N2").ClearContents()
$ws.Range("H5").Value = 61.444443
$ws.Range("I5").Value = 50.42857
$ws.Range("K5").Value = 50.42857
$ws.Range("M5").Value = 61.57143
$ws.Range("N5").ClearContents()
$ws.Range("H61").Value = 5880.353
$ws.Range("I61").Value = 4059
$ws.Range("K61").Value = 4059
$ws.Range("M61").Value = -3847
$ws.Range("N61").ClearContents()
$ws.Range("H74").Value = 2398.1052
$ws.Range("I74").Value = 1975.7778
$ws.Range("K74").Value = 1975.7778
$ws.Range("M74").Value = -1101.7778
$ws.Range("N74").ClearContents()
$ws.Range("H77").Value = 2398.1052
$ws.Range("I77").Value = 1975.7778
$ws.Range("K77").Value = 9878.889000000001
$ws.Range("M77").Value = -5510.889000000001
$ws.Range("N77").ClearContents()
$ws.Range("H116").Value = 740.0769
$ws.Range("I116").Value = 685.0833
$ws.Range("K116").Value = 685.0833
$ws.Range("M116").Value = 1608.9167
$ws.Range("N116").ClearContents()
$ws.Range("H132").Value = 2947.4443
$ws.Range("I132").Value = 2563.647
$ws.Range("J132").Value = 3599.9
$ws.Range("K132").Value = 7690.941
$ws.Range("L132").Value = 10799.7
$ws.Range("M132").Value = -5160.941
$ws.Range("N132").Value = -15859.7
$ws.Range("H136").Value = 5880.353
$ws.Range("I136").Value = 4059
$ws.Range("K136").Value = 12177
$ws.Range("M136").Value = -9627
$ws.Range("N136").ClearContents()
$ws.Range("H140").Value = 0
$ws.Range("J140").Value = 0
$ws.Range("L140").Value = 0
$ws.Range("N140").ClearContents()
$ws.Range("H141").Value = 0
$ws.Range("J141").Value = 0
$ws.Range("L141").Value = 0
$ws.Range("N141").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 740.0769
$ws.Range("I3").Value = 685.0833
$ws.Range("K3").Value = 685.0833
$ws.Range("M3").Value = -571.0833
$ws.Range("N3").ClearContents()
$ws.Range("H4").Value = 61.444443
$ws.Range("I4").Value = 50.42857
$ws.Range("K4").Value = 50.42857
$ws.Range("M4").Value = 64.57142999999999
$ws.Range("N4").ClearContents()
$ws.Range("H22").Value = 1458.1
$ws.Range("I22").Value = 1458.1
$ws.Range("K22").Value = 1458.1
$ws.Range("M22").Value = -1285.1

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 4001.3333
$ws.Range("I16").Value = 3995
$ws.Range("K16").Value = 3995
$ws.Range("M16").Value = -3708
$ws.Range("H31").Value = 4743.3887
$ws.Range("I31").Value = 2186.4
$ws.Range("J31").Value = 7939.625
$ws.Range("K31").Value = 2186.4
$ws.Range("L31").Value = 7939.625
$ws.Range("M31").Value = -1891.4
$ws.Range("N31").Value = -8529.625
$ws.Range("H34").Value = 4743.3887
$ws.Range("I34").Value = 2186.4
$ws.Range("J34").Value = 7939.625
$ws.Range("K34").Value = 2186.4
$ws.Range("L34").Value = 7939.625
$ws.Range("M34").Value = -1984.4
$ws.Range("N34").Value = -8343.625
$ws.Range("H107").Value = 271.33334
$ws.Range("I107").Value = 271.33334
$ws.Range("K107").Value = 271.33334
$ws.Range("M107").Value = 1648.66666
$ws.Range("H113").Value = 4001.3333
$ws.Range("I113").Value = 3995
$ws.Range("K113").Value = 3995
$ws.Range("M113").Value = -1825
$ws.Range("H140").Value = 0
$ws.Range("J140").Value = 0
$ws.Range("L140").Value = 0
$ws.Range("N140").ClearContents()
$ws.Range("H141").Value = 69000
$ws.Range("J141").Value = 69000
$ws.Range("L141").Value = 69000
$ws.Range("N141").Value = -79360

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 102.55556
$ws.Range("I12").Value = 10.166667
$ws.Range("J12").Value = 148.75
$ws.Range("K12").Value = 30.500001
$ws.Range("L12").Value = 446.25
$ws.Range("M12").Value = 142.499999
$ws.Range("N12").Value = -792.25
$ws.Range("H23").Value = 219.4
$ws.Range("J23").Value = 250
$ws.Range("L23").Value = 750
$ws.Range("N23").Value = -1220
$ws.Range("H129").Value = 1909
$ws.Range("J129").Value = 1895
$ws.Range("L129").Value = 5685
$ws.Range("N129").Value = -15685
$ws.Range("H134").Value = 4076.6667
$ws.Range("J134").Value = 10000
$ws.Range("L134").Value = 30000
$ws.Range("N134").Value = -40140

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3499.5
$ws.Range("I80").Value = 3499
$ws.Range("K80").Value = 3499
$ws.Range("M80").Value = -2501
$ws.Range("H83").Value = 3499.5
$ws.Range("I83").Value = 3499
$ws.Range("K83").Value = 17495
$ws.Range("M83").Value = -12503
$ws.Range("H113").Value = 6910.3335
$ws.Range("I113").Value = 5529.4287
$ws.Range("K113").Value = 5529.4287
$ws.Range("M113").Value = -3359.4287
$ws.Range("N113").ClearContents()
$ws.Range("H132").Value = 1581.3
$ws.Range("I132").Value = 1112.4706
$ws.Range("K132").Value = 3337.4118
$ws.Range("M132").Value = -807.4118000000003
$ws.Range("N132").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 1000
$ws.Range("I93").Value = 1000
$ws.Range("K93").Value = 1000
$ws.Range("M93").Value = 248
$ws.Range("H100").Value = 6303.769
$ws.Range("I100").Value = 3993.625
$ws.Range("K100").Value = 3993.625
$ws.Range("M100").Value = -3452.625
$ws.Range("N100").ClearContents()
$ws.Range("H136").Value = 1998
$ws.Range("I136").Value = 1997.6
$ws.Range("K136").Value = 5992.799999999999
$ws.Range("M136").Value = -3442.799999999999
$ws.Range("N136").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H45").Value = 21996.334
$ws.Range("I45").Value = 17999.5
$ws.Range("J45").Value = 29990
$ws.Range("K45").Value = 17999.5
$ws.Range("L45").Value = 29990
$ws.Range("M45").Value = -17508.5
$ws.Range("N45").Value = -30972
$ws.Range("H49").Value = 30000
$ws.Range("J49").Value = 30000
$ws.Range("L49").Value = 30000
$ws.Range("N49").Value = -30460
